$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Acompanhar e" + "  " + "avaliar o desempenho..." -> merge into a
#    single run (drops the gramStart/gramEnd proofErr wrapping the
#    double space introduced by the spell/grammar checker).
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Acompanhar e  avaliar o desempenho da equipe,propondo melhorias no desenvolvimento do aplicativo.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Acompanhar e  avaliar o desempenho da equipe,propondo melhorias no desenvolvimento do aplicativo.",
    2) | Out-Null

# ---------------------------------------------------------------------
# Table 4 ("Necessidade / Prioridade / Característica / Release
# Planejado") — fill in the previously-empty rows.
# ---------------------------------------------------------------------
$t4 = $d.Tables.Item(4)

# Row 6 — "Buscar"
$rowBuscar = $t4.Rows.Item(6)
$c1 = $rowBuscar.Cells.Item(1).Range
$c1.End = $c1.End - 1
$c1.Collapse(0)
$c1.InsertAfter(" informações na Base de dados")
$rowBuscar.Cells.Item(3).Range.Text = "Realizar pesquisa na base de dados e exibi na tela para o usuário."
$rowBuscar.Cells.Item(4).Range.Text = "E2"

# Row 7 — "Gerar relatorio"
$rowRelatorio = $t4.Rows.Item(7)
$rowRelatorio.Cells.Item(3).Range.Text = "Gera um relatório contendo o balanço de coletas do mês."
$rowRelatorio.Cells.Item(4).Range.Text = "E2"

# Row 8 — "Autenticação"
$rowAuth = $t4.Rows.Item(8)
$cAuth = $rowAuth.Cells.Item(1).Range
$cAuth.End = $cAuth.End - 1
$cAuth.Collapse(0)
$cAuth.InsertAfter(" de usuário ")
$rowAuth.Cells.Item(3).Range.Text = "Autenticação do usuário para utilizar serviços do aplicativo."
$rowAuth.Cells.Item(4).Range.Text = "E1"

# ---------------------------------------------------------------------
# Table 5 ("Requisito / Prioridade / Release Planejado") — rename the
# last row and append new "Confiabilidade" / "Desempenho" rows.
# ---------------------------------------------------------------------
$t5 = $d.Tables.Item(5)

$d.Content.Find.Execute(
    "Requisitos não funcionais", $true, $false, $false, $false, $false,
    $true, 1, $false, "Interface intuitiva", 2) | Out-Null

$rowInterface = $t5.Rows.Item(3)
$rowInterface.Cells.Item(2).Range.Text = "Média"

$rowConf = $t5.Rows.Add()
$rowConf.Cells.Item(1).Range.Text = "Confiabilidade "
$rowConf.Cells.Item(2).Range.Text = "Baixa"

$rowDesemp = $t5.Rows.Add()
$rowDesemp.Cells.Item(1).Range.Text = "Desempenho"
$rowDesemp.Cells.Item(2).Range.Text = "Alta"
